$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 83; existing rows 83-85 shift down to 84-86
$ws.Rows("83:83").Insert()

# Populate the new row 83 with the latest weekly price entry
$ws.Range("A83").Value2 = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value2 = 44448
$ws.Range("E83").Value2 = 8
$ws.Range("F83").Value2 = 100112003
$ws.Range("G83").Value = "Ajo"
$ws.Range("H83").Value = "Chino"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value2 = 400
$ws.Range("K83").Value2 = 14000
$ws.Range("L83").Value2 = 14500
$ws.Range("M83").Value2 = 14250
$ws.Range("N83").Value = "$/caja 10 kilos"
$ws.Range("O83").Value = "China"
$ws.Range("P83").Value2 = 1425
$ws.Range("Q83").Value2 = 10
$ws.Range("R83").Value = "Hortaliza"
